# Weekly price-sheet refresh: a new observation is inserted at row 47
# (pushing every subsequent row down by one, with the former last row
# re-appearing as the new last row), and the freshly opened row 47 is
# filled with the new week's Fecha (date serial) and Volumen values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 47; rows 47..238 shift down to 48..239
# and the sheet grows from A1:R238 to A1:R239.
$ws.Rows.Item(47).Insert()

# Seed the new row 47 with the same "static" field values (market,
# region, category, quality, prices, units, origin, classification...)
# as the row right below it (the row that used to be row 47), then
# overwrite just the date and volume with this week's new figures.
$ws.Rows.Item(48).Copy()
$ws.Rows.Item(47).PasteSpecial()

$ws.Range("D47").Value = 44565
$ws.Range("J47").Value = 180
